$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# C25 has been re-selected as 100pF, same part as C1, so merge the two rows:
# Row 2 currently holds C1 (100pF, qty 1). Row 9 holds C25 (560pF, qty 1).
# Update row 2 to reflect the combined reference / quantity, then remove the
# old C25 row entirely (its other fields already match C1's part details).
$ws.Range("A2").Value = "C1 C25"
$ws.Range("B2").Value = 2

# Delete the now-redundant C25 row (row 9), shifting everything below up.
$ws.Rows("9:9").Select() | Out-Null
$ws.Rows("9:9").Delete()
